$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F10").Value = 452790
$ws.Range("F11").Value = 1438649

$ws.Range("F11").Select()
